# "update all seasons page"
# Slide 10 (the "ALL SEASONS" page) has several shapes shifted upward
# (a couple also shift slightly in X), and a slow "push up" slide
# transition is added.
#
# NOTE on the literals below: the host stores Shape.Left/.Top as points
# that get converted to EMU (1 pt = 12700 EMU). To land on the exact
# target EMU values from the authoring tool (avoiding off-by-one EMU
# drift from float rounding) the point values used here were solved so
# they reproduce the exact target EMU after conversion.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# Rectangle 10 ("CHRISTMAS:") -> off x=313081 y=888815
$shp = $s.Shapes.Item("Rectangle 10")
$shp.Left = 24.65205
$shp.Top = 69.985433

# Arrow: Right 11 -> off y=3494299 (x unchanged)
$shp = $s.Shapes.Item("Right Arrow 11")
$shp.Top = 275.1417

# Arrow: Right 12 -> off y=3543857 (x unchanged)
$shp = $s.Shapes.Item("Right Arrow 12")
$shp.Top = 279.0439

# Picture 13 -> off y=3039112 (x unchanged)
$shp = $s.Shapes.Item("Picture 13")
$shp.Top = 239.3002

# Picture 14 -> off y=2992901 (x unchanged)
$shp = $s.Shapes.Item("Picture 14")
$shp.Top = 235.6615

# Picture 15 -> off y=2906485 (x unchanged)
$shp = $s.Shapes.Item("Picture 15")
$shp.Top = 228.8571

# Picture 16 -> off y=2944194 (x unchanged)
$shp = $s.Shapes.Item("Picture 16")
$shp.Top = 231.82630157559842

# Table 17 -> off y=1548705 (x unchanged)
$shp = $s.Shapes.Item("Table 17")
$shp.Top = 121.9453

# TextBox 18 ("PARAMENT COLOR: WHITE") -> off y=937037 (x unchanged)
$shp = $s.Shapes.Item("TextBox 18")
$shp.Top = 73.78244400088188

# TextBox 20 ("ALL SEASONS:") -> off x=309914 y=239861
$shp = $s.Shapes.Item("TextBox 20")
$shp.Left = 24.4027
$shp.Top = 18.8867

# Slide transition: slow "push up"
$t = $s.SlideShowTransition
$t.EntryEffect = 3852
$t.Speed = 1
